# Apply the commit: rename "45 to 94 Years" sheet/labels to "45 to 54 Years",
# and normalize placeholder "0" text cells to "0.0" text cells across a few
# age-bracket sheets (troubleshooting clustered bar chart source data).

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "45 to 94 Years" sheet to "45 to 54 Years" -------------
$wsOld = $wb.Worksheets.Item("45 to 94 Years")
$wsOld.Name = "45 to 54 Years"

# --- 2. Update the column header label on that sheet (B1) ------------------
$wsRenamed = $wb.Worksheets.Item("45 to 54 Years")
$wsRenamed.Range("B1").Value = "General Warehousing & Storage - 45 to 54 Years"

# --- 3. Update the matching header label on the "df_all" sheet (K1) --------
$wsAll = $wb.Worksheets.Item("df_all")
$wsAll.Range("K1").Value = "General Warehousing & Storage - 45 to 54 Years"

# Helper: write a string value into a cell while keeping it text (not
# auto-coerced to a number) and without leaving a lingering style/number
# format behind on the cell.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- 4. Rows 2-12 of column E & F on "df_all": "0" -> "0.0" -----------------
for ($r = 2; $r -le 12; $r++) {
    Set-TextValue $wsAll.Range("E$r") "0.0"
    Set-TextValue $wsAll.Range("F$r") "0.0"
}

# --- 5. Rows 2-12 of column B on "Under 16 Years": "0" -> "0.0" -------------
$wsUnder16 = $wb.Worksheets.Item("Under 16 Years")
for ($r = 2; $r -le 12; $r++) {
    Set-TextValue $wsUnder16.Range("B$r") "0.0"
}

# --- 6. Rows 2-12 of column B on "16 to 17 Years": "0" -> "0.0" -------------
$ws1617 = $wb.Worksheets.Item("16 to 17 Years")
for ($r = 2; $r -le 12; $r++) {
    Set-TextValue $ws1617.Range("B$r") "0.0"
}
